# Update cryptocurrency price/volume-label data on Sheet1 (commit: "Updated symbol list on
# Tue Dec 27 10:50:39 UTC 2022 with GitHub Actions").
#
# The Price column (D) stores values as text (e.g. "243.62"), mirroring the source sheet's
# inline-string cells. Plain `.Value = "<number-looking text>"` assignment would let Excel
# auto-coerce the text to a real number, so each D-column write briefly forces the cell to
# Text format, assigns the literal string, then restores the cell to the default "Normal"
# style so no stray number-format/style is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue "D2" "243.54"
Set-TextValue "D4" "5.393"
Set-TextValue "D5" "0.05954"
Set-TextValue "D6" "3.435"
Set-TextValue "D7" "6.496"
Set-TextValue "D8" "0.8103"
Set-TextValue "D9" "0.9306"
Set-TextValue "D10" "0.1431"
Set-TextValue "D11" "0.07408"
Set-TextValue "D13" "0.03079"
Set-TextValue "D15" "3.861"
Set-TextValue "D16" "0.001575"
Set-TextValue "D17" "0.04697"
Set-TextValue "D18" "0.0005981"
$ws.Range("E18").Value = "17OneONE"
Set-TextValue "D19" "0.005966"
Set-TextValue "D20" "0.001261"
$ws.Range("E20").Value = "19BitKanKAN"
Set-TextValue "D21" "0.004787"
Set-TextValue "D22" "0.00008002"
Set-TextValue "D23" "3.571"
Set-TextValue "D24" "2.134"
Set-TextValue "D27" "0.0002340"
Set-TextValue "D40" "0.03938"
Set-TextValue "D41" "0.006321"
Set-TextValue "D43" "0.003501"
$ws.Range("E43").Value = "42CEJICEJIBestin24h"
Set-TextValue "D44" "0.008961"
Set-TextValue "D45" "0.00005190"
Set-TextValue "D47" "0.6851"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
Set-TextValue "D48" "0.002070"
Set-TextValue "D49" "0.00002100"
Set-TextValue "D50" "0.0002000"
